$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

# --- Update the existing "Day 2" block totals ---
$ws.Range("C9").Value = 78
$ws.Range("C11").Value = 114

# --- Add the new "Day 3" block (rows 13-16), cloned from the "Day 2" block (rows 8-11) ---
# Copying the formatted range preserves cell styles/merge so the new block
# looks identical to the Day 2 block it is based on.
$ws.Range("B8:C11").Copy($ws.Range("B13"))

# Give the new block its own heading text and totals
$ws.Range("B13").Value = "Spint( 35) - Day 3 - Test Case Summary"
$ws.Range("C14").Value = 117
$ws.Range("C15").Value = 75
$ws.Range("C16").Value = 154

# Match the row height used by the rest of the summary blocks
$ws.Range("B13:C16").RowHeight = 18

# --- Update the recorded selection ---
$ws.Range("G15").Select()
